# Applies the bilibili-scraped Suzhou/Kunshan/Changshu listing refresh
# (commit 456a3b4) to the worksheets that carry the event table: '展览'
# (sheet1) and '全部类型' (sheet4) hold identical tables in this workbook.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # -- counts refreshed in place (no row/col shift) --
    $ws.Cells.Item(2, 6).Value = 848
    $ws.Cells.Item(4, 6).Value = 2162
    $ws.Cells.Item(6, 6).Value = 12566
    $ws.Cells.Item(8, 6).Value = 106
    $ws.Cells.Item(10, 6).Value = 457
    $ws.Cells.Item(12, 6).Value = 942
    $ws.Cells.Item(13, 6).Value = 13650
    $ws.Cells.Item(14, 6).Value = 13964

    # -- Rows 19-29: five new events were added (2024-05-02 meet-and-greets
    #    at the Suzhou Bay theatre), pushing the previous rows 19-24 down to
    #    25-29 (two of which also had their own field tweaks). Rewrite every
    #    field of rows 19-29 to the final published content. --

    # Row 19
    $ws.Cells.Item(19, 1).Value = 18
    $ws.Cells.Item(19, 2).Value = "'2024-05-02"
    $ws.Cells.Item(19, 2).Style = "Normal"
    $ws.Cells.Item(19, 3).Value = "苏州·动漫游戏嘉年华"
    $ws.Cells.Item(19, 4).Value = "东太湖大道12000号 苏州湾大剧院"
    $ws.Cells.Item(19, 5).Value = "2024.05.02 10:00-05.02 17:00"
    $ws.Cells.Item(19, 6).Value = 5
    $ws.Cells.Item(19, 7).Value = 58
    $ws.Cells.Item(19, 8).Value = "https://show.bilibili.com/platform/detail.html?id=82824"
    $ws.Cells.Item(19, 9).Value = "//i1.hdslb.com/bfs/openplatform/202403/HzWBEJeE1710324788092.jpeg"

    # Row 20
    $ws.Cells.Item(20, 1).Value = 19
    $ws.Cells.Item(20, 2).Value = "'2024-05-02"
    $ws.Cells.Item(20, 2).Style = "Normal"
    $ws.Cells.Item(20, 3).Value = "苏州·动漫游戏嘉年华cv刘圣博见面会"
    $ws.Cells.Item(20, 4).Value = "东太湖大道12000号 苏州湾大剧院"
    $ws.Cells.Item(20, 5).Value = "2024.05.02 10:00-05.02 17:00"
    $ws.Cells.Item(20, 6).Value = 1
    $ws.Cells.Item(20, 7).Value = 188
    $ws.Cells.Item(20, 8).Value = "https://show.bilibili.com/platform/detail.html?id=83038"
    $ws.Cells.Item(20, 9).Value = "//i0.hdslb.com/bfs/openplatform/202403/D94B39u21710901393375.jpeg"

    # Row 21
    $ws.Cells.Item(21, 1).Value = 20
    $ws.Cells.Item(21, 2).Value = "'2024-05-02"
    $ws.Cells.Item(21, 2).Style = "Normal"
    $ws.Cells.Item(21, 3).Value = "苏州·动漫游戏嘉年华cv张文钰见面会"
    $ws.Cells.Item(21, 4).Value = "东太湖大道12000号 苏州湾大剧院"
    $ws.Cells.Item(21, 5).Value = "2024.05.02 10:00-05.02 17:00"
    $ws.Cells.Item(21, 6).Value = 2
    $ws.Cells.Item(21, 7).Value = 188
    $ws.Cells.Item(21, 8).Value = "https://show.bilibili.com/platform/detail.html?id=83037"
    $ws.Cells.Item(21, 9).Value = "//i1.hdslb.com/bfs/openplatform/202403/eP5thEuS1710901472569.jpeg"

    # Row 22
    $ws.Cells.Item(22, 1).Value = 21
    $ws.Cells.Item(22, 2).Value = "'2024-05-02"
    $ws.Cells.Item(22, 2).Style = "Normal"
    $ws.Cells.Item(22, 3).Value = "苏州·动漫游戏嘉年华cv沐霏见面会"
    $ws.Cells.Item(22, 4).Value = "东太湖大道12000号 苏州湾大剧院"
    $ws.Cells.Item(22, 5).Value = "2024.05.02 10:00-05.02 17:00"
    $ws.Cells.Item(22, 6).Value = 4
    $ws.Cells.Item(22, 7).Value = 188
    $ws.Cells.Item(22, 8).Value = "https://show.bilibili.com/platform/detail.html?id=82891"
    $ws.Cells.Item(22, 9).Value = "//i1.hdslb.com/bfs/openplatform/202403/8VORpvQc1710900704258.jpeg"

    # Row 23
    $ws.Cells.Item(23, 1).Value = 22
    $ws.Cells.Item(23, 2).Value = "'2024-05-03"
    $ws.Cells.Item(23, 2).Style = "Normal"
    $ws.Cells.Item(23, 3).Value = "常熟·CDW·动漫展03"
    $ws.Cells.Item(23, 4).Value = "常熟国际展览中心 国际展览中心"
    $ws.Cells.Item(23, 5).Value = "2024.05.03 09:00-05.04 17:30"
    $ws.Cells.Item(23, 6).Value = 1049
    $ws.Cells.Item(23, 7).Value = 60
    $ws.Cells.Item(23, 8).Value = "https://show.bilibili.com/platform/detail.html?id=82489"
    $ws.Cells.Item(23, 9).Value = "//i0.hdslb.com/bfs/openplatform/202403/XK411blC1709794808211.jpeg"

    # Row 24
    $ws.Cells.Item(24, 1).Value = 23
    $ws.Cells.Item(24, 2).Value = "'2024-05-03"
    $ws.Cells.Item(24, 2).Style = "Normal"
    $ws.Cells.Item(24, 3).Value = "昆山·第十二届理想乡动漫游戏展嘉宾沈辞签售会"
    $ws.Cells.Item(24, 4).Value = "花桥经济开发区绿地大道1598号 花桥国际博览中心"
    $ws.Cells.Item(24, 5).Value = "2024.05.03 14:00-05.03 16:00"
    $ws.Cells.Item(24, 6).Value = 106
    $ws.Cells.Item(24, 7).Value = 1
    $ws.Cells.Item(24, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81120"
    $ws.Cells.Item(24, 9).Value = "//i0.hdslb.com/bfs/openplatform/202401/4Pay1rR61705648901961.jpeg"

    # Row 25
    $ws.Cells.Item(25, 1).Value = 24
    $ws.Cells.Item(25, 2).Value = "'2024-05-03"
    $ws.Cells.Item(25, 2).Style = "Normal"
    $ws.Cells.Item(25, 3).Value = "昆山·第十二届理想乡动漫游戏展嘉宾矮乐多aliga签售会"
    $ws.Cells.Item(25, 4).Value = "花桥经济开发区绿地大道1598号 花桥国际博览中心"
    $ws.Cells.Item(25, 5).Value = "2024.05.03 14:00-05.03 16:00"
    $ws.Cells.Item(25, 6).Value = 54
    $ws.Cells.Item(25, 7).Value = 1
    $ws.Cells.Item(25, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81114"
    $ws.Cells.Item(25, 9).Value = "//i1.hdslb.com/bfs/openplatform/202401/Peub7FOc1705648580577.jpeg"

    # Row 26
    $ws.Cells.Item(26, 1).Value = 25
    $ws.Cells.Item(26, 2).Value = "'2024-05-04"
    $ws.Cells.Item(26, 2).Style = "Normal"
    $ws.Cells.Item(26, 3).Value = "【大会员提前抢】苏州·OCG国潮动漫游戏嘉年华阿杰内场"
    $ws.Cells.Item(26, 4).Value = "苏州大道东688号 苏州国际博览中心"
    $ws.Cells.Item(26, 5).Value = "2024.05.04 09:00-05.04 17:00"
    $ws.Cells.Item(26, 6).Value = 520
    $ws.Cells.Item(26, 7).Value = 288
    $ws.Cells.Item(26, 8).Value = "https://show.bilibili.com/platform/detail.html?id=82940"
    $ws.Cells.Item(26, 9).Value = "//i2.hdslb.com/bfs/openplatform/202403/lLKmv48C1710511298160.jpeg"

    # Row 27
    $ws.Cells.Item(27, 1).Value = 26
    $ws.Cells.Item(27, 2).Value = "'2024-05-04"
    $ws.Cells.Item(27, 2).Style = "Normal"
    $ws.Cells.Item(27, 3).Value = "苏州·OCG国潮动漫游戏嘉年华"
    $ws.Cells.Item(27, 4).Value = "苏州大道东688号 苏州国际博览中心"
    $ws.Cells.Item(27, 5).Value = "2024.05.04 09:00-05.05 17:00"
    $ws.Cells.Item(27, 6).Value = 5067
    $ws.Cells.Item(27, 7).Value = 65
    $ws.Cells.Item(27, 8).Value = "https://show.bilibili.com/platform/detail.html?id=82779"
    $ws.Cells.Item(27, 9).Value = "//i1.hdslb.com/bfs/openplatform/202403/hcgdIzw61710298907237.jpeg"

    # Row 28
    $ws.Cells.Item(28, 1).Value = 27
    $ws.Cells.Item(28, 2).Value = "'2024-05-18"
    $ws.Cells.Item(28, 2).Style = "Normal"
    $ws.Cells.Item(28, 3).Value = "苏州·YoungComic动漫嘉年华"
    $ws.Cells.Item(28, 4).Value = "清禾路886号 尹山湖大剧院"
    $ws.Cells.Item(28, 5).Value = "2024.05.18 10:00-05.18 17:00"
    $ws.Cells.Item(28, 6).Value = 0
    $ws.Cells.Item(28, 7).Value = 60
    $ws.Cells.Item(28, 8).Value = "https://show.bilibili.com/platform/detail.html?id=83142"
    $ws.Cells.Item(28, 9).Value = "//i2.hdslb.com/bfs/openplatform/202403/4wWLK6Jg1710840463319.jpeg"

    # Row 29
    $ws.Cells.Item(29, 1).Value = 28
    $ws.Cells.Item(29, 2).Value = "'2024-06-08"
    $ws.Cells.Item(29, 2).Style = "Normal"
    $ws.Cells.Item(29, 3).Value = "【会员购严选】苏州·Come in joy动漫国潮文化节"
    $ws.Cells.Item(29, 4).Value = "金山南路288号 广电国际会展中心"
    $ws.Cells.Item(29, 5).Value = "2024.06.08 10:00-06.09 17:00"
    $ws.Cells.Item(29, 6).Value = 250
    $ws.Cells.Item(29, 7).Value = 60
    $ws.Cells.Item(29, 8).Value = "https://show.bilibili.com/platform/detail.html?id=82233"
    $ws.Cells.Item(29, 9).Value = "//i0.hdslb.com/bfs/openplatform/202403/F86lgbSt1709278264141.jpeg"

    # Rows 25-29 are brand-new cells in column A; copy the bordered/bold/
    # centered style used by the existing index column (e.g. A24) onto them
    # so the new rows match the look of the rest of the column, then
    # restore their numeric index values (PasteSpecial also pastes the
    # source's value, which we don't want to keep).
    $ws.Range("A24").Copy()
    $ws.Range("A25:A29").PasteSpecial(-4122)
    $ws.Cells.Item(25, 1).Value = 24
    $ws.Cells.Item(26, 1).Value = 25
    $ws.Cells.Item(27, 1).Value = 26
    $ws.Cells.Item(28, 1).Value = 27
    $ws.Cells.Item(29, 1).Value = 28
    $excel.CutCopyMode = 0
}

